# "Rescaled the game by 0.5" - halve every weapon-effect / projectile-size
# distance constant in column N. Column M either holds a literal 0 (left
# alone - the halving doesn't change it), a formula =N/0.3*0.25 (recalculates
# automatically once N changes), or, in the three summary rows at the bottom
# (60-62), a literal value that also needs to be halved explicitly.
# Also update the sheet's scroll position / selection to match the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warheads")

# Rows whose column-N "distance" constant must be halved.
$nRows = @(2,3,4,6,7,8,10,11,12,20,21,22,24,25,26,28,29,30,38,39,40,42,43,44,46,47,48,56,57,60,61,62)

foreach ($r in $nRows) {
    $cell = $ws.Cells.Item($r, 14)   # column N
    $cell.Value2 = $cell.Value2 / 2
}

# The three bottom summary rows also carry a literal (non-formula) constant
# in column M that needs the same treatment.
$mRows = @(60,61,62)
foreach ($r in $mRows) {
    $cell = $ws.Cells.Item($r, 13)   # column M
    $cell.Value2 = $cell.Value2 / 2
}

# Match the saved view: scrolled down a bit further, whole table selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 44
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1:N62").Select()
